$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Unprotect("D382")

$ws.Range("A38").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-25 for illustrative purposes only and are subject to change."

$ws.Range("D2").Value = 0.0353464284896958
$ws.Range("D3").Value = 0.02024104734323187
$ws.Range("E3").Value = 0.001173249902229223
$ws.Range("D4").Value = 0.01926698131834123
$ws.Range("E4").Value = 0.0008006405124099114
$ws.Range("D5").Value = 0.03786455864624916
$ws.Range("E5").Value = 0
$ws.Range("D6").Value = 0.03411113658982407
$ws.Range("E6").Value = 0.0004001600640257674
$ws.Range("D7").Value = 0.01970611277173796
$ws.Range("E7").Value = 0.0003862495171880198
$ws.Range("D8").Value = 0.03726528172637219
$ws.Range("E8").Value = -0.01118615449556781
$ws.Range("D9").Value = 0.02033360289800373
$ws.Range("E9").Value = 0.004851316144102169
$ws.Range("D10").Value = 0.02595767217299486
$ws.Range("E10").Value = -0.009774215619196558
$ws.Range("D11").Value = 0.02409666615805089
$ws.Range("E11").Value = -0.006843906291129298
$ws.Range("D12").Value = 0.05748598105564379
$ws.Range("E12").Value = -0.005420692905962587
$ws.Range("D13").Value = 0.02477078926394129
$ws.Range("E13").Value = 0.002949852507374562
$ws.Range("D14").Value = 0.02680477876913595
$ws.Range("E14").Value = -0.005277864017386058
$ws.Range("D15").Value = 0.03238307769851072
$ws.Range("E15").Value = -0.008775008775008808
$ws.Range("D16").Value = 0.01959483298466851
$ws.Range("E16").Value = -0.01205727204220053
$ws.Range("D17").Value = 0.031620610721536
$ws.Range("E17").Value = 0.005375912701596652
$ws.Range("D18").Value = 0.04195851816324111
$ws.Range("E18").Value = -0.0009191176470590978
$ws.Range("D19").Value = 0.125397046331007
$ws.Range("E19").Value = -0.001328903654485236
$ws.Range("D20").Value = 0.009223536532524879
$ws.Range("E20").Value = -0.01935973680880676
$ws.Range("D21").Value = 0.01537137688007911
$ws.Range("E21").Value = -0.01247833622183714
$ws.Range("D22").Value = 0.01754597572862508
$ws.Range("E22").Value = 0.004916420845624492
$ws.Range("D23").Value = 0.01544810056363998
$ws.Range("E23").Value = 0.008592910848550073
$ws.Range("D24").Value = 0.02169301260871535
$ws.Range("E24").Value = 0.0009052504526252392
$ws.Range("D25").Value = 0.01277982134646716
$ws.Range("E25").Value = -0.01729183293429115
$ws.Range("D26").Value = 0.04242180336886996
$ws.Range("E26").Value = -0.0005920981806438741
$ws.Range("D27").Value = 0.02381103281495288
$ws.Range("E27").Value = -0.0000980296049406526
$ws.Range("D28").Value = 0.04561101204069623
$ws.Range("E28").Value = -0.003317535545023675
$ws.Range("D29").Value = 0.05610348319941616
$ws.Range("E29").Value = -0.002481829462861129
$ws.Range("D30").Value = 0.01301410259448341
$ws.Range("E30").Value = -0.007766990291261933
$ws.Range("D31").Value = 0.02054388262331873
$ws.Range("E31").Value = -0.00191424196018386
$ws.Range("D32").Value = 0.01351671376225492
$ws.Range("E32").Value = 0
$ws.Range("D33").Value = 0.04161940151162909
$ws.Range("E33").Value = -0.001031459515213884
$ws.Range("D34").Value = 0.01709162132214087
$ws.Range("E34").Value = 0.0004364271166716893
$ws.Range("E35").Value = -0.002458659153158238

$ws.Protect("D382")
